$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.260.88'
$ws.Range('E2').Value = '  +3.43%  '
$ws.Range('D3').Value = '2.367.58'
$ws.Range('E3').Value = '  +1.88%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '543.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +1.01%  '
$ws.Range('D9').Value = '2.368.36'
$ws.Range('E9').Value = '  +1.93%  '
$ws.Range('E10').Value = '  +2.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.44'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.18%  '
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('E13').Value = '  +5.85%  '
$ws.Range('D14').Value = '2.791.56'
$ws.Range('E14').Value = '  +1.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.65'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').Value = '58.201.24'
$ws.Range('E16').Value = '  +3.35%  '
$ws.Range('E17').Value = '  +1.66%  '
$ws.Range('D18').Value = '2.357.89'
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '339.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.23'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.86'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.21%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.76%  '
$ws.Range('E25').Value = '  +4.03%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('E28').Value = '  +9.17%  '
$ws.Range('E29').Value = '  +6.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.98'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.94%  '
$ws.Range('D31').Value = '0.0₃0737'
$ws.Range('E31').Value = '  +3.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.58'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.92%  '
$ws.Range('E34').Value = '  +16.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.16'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.14%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.26'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.51%  '
$ws.Range('E39').Value = '  +5.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '39.47'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '149.77'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.379'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.64'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '284.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.35'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.70%  '
$ws.Range('E46').Value = '  +0.98%  '
$ws.Range('E47').Value = '  +2.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.561'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('E49').Value = '  +2.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.57'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.381'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.20%  '
